# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets.
#
# Mapping of row -> (old value, new value) is identical between the two
# sheets (the "全部类型" sheet has one extra row inserted near the top,
# which shifts everything below it down by one row relative to "展览").

$wb = $excel.ActiveWorkbook

$updates1 = @{
    2  = 1906
    3  = 513
    5  = 182
    6  = 2755
    10 = 1579
    11 = 558
    13 = 342
    14 = 234
    22 = 14
    23 = 224
    25 = 1756
    27 = 425
    28 = 86
    31 = 312
    32 = 453
}

$updates4 = @{
    2  = 1906
    4  = 513
    6  = 182
    7  = 2755
    11 = 1579
    12 = 558
    14 = 342
    15 = 234
    23 = 14
    24 = 224
    26 = 1756
    28 = 425
    29 = 86
    32 = 312
    33 = 453
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
